$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.587.40"
$ws.Range("E2").Value = "  +0.25%  "

$ws.Range("E3").Value = "  +0.42%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.88%  "

$ws.Range("E7").Value = "  -0.53%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.542"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.49"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.78%  "

$ws.Range("E11").Value = "  -0.19%  "

$ws.Range("E12").Value = "  +0.68%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.62"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.26%  "

$ws.Range("E14").Value = "  +0.87%  "

$ws.Range("D15").Value = "2.885.01"
$ws.Range("E15").Value = "  +0.30%  "

$ws.Range("D16").Value = "2.496.23"
$ws.Range("E16").Value = "  +0.29%  "

$ws.Range("E17").Value = "  +0.96%  "

$ws.Range("D18").Value = "47.489.84"
$ws.Range("E18").Value = "  +0.22%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.93%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.64"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.41%  "

$ws.Range("D21").Value = "0.0₃0942"
$ws.Range("E21").Value = "  +0.79%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +15.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.68"
$ws.Range("D23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "246.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.43%  "

$ws.Range("E25").Value = "  -0.56%  "

$ws.Range("E26").Value = "  +0.08%  "

$ws.Range("E27").Value = "  -1.08%  "

$ws.Range("E28").Value = "  +0.08%  "

$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.140"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.81%  "

$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.20%  "

$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.08"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.54%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.92"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.01%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.37"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.77%  "

$ws.Range("E34").Value = "  -0.25%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0789"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.03%  "

$ws.Range("E36").Value = "  +0.17%  "

$ws.Range("E37").Value = "  +3.01%  "

$ws.Range("E38").Value = "  +1.09%  "

$ws.Range("E39").Value = "  -1.05%  "

$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "22.62"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.42%  "

$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.112"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.46%  "

$ws.Range("E42").Value = "  -1.87%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "119.12"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.74%  "

$ws.Range("E44").Value = "  +0.25%  "

$ws.Range("D45").Value = "1.999.12"
$ws.Range("E45").Value = "  +1.89%  "

$ws.Range("E46").Value = "  +2.37%  "

$ws.Range("E47").Value = "  -2.74%  "

$ws.Range("E48").Value = "  +0.19%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.06"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.58%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.09%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "56.69"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.63%  "
